# Insert a new data row right after the existing "2026/02/26" (木) block,
# shifting the remainder of the "2026/12/29 ..." onward data down by one
# row, matching a new entry: 2026/02/26, 木, 6, 201 at row 876.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 876 and everything below it down by one row to make room for
# the newly-inserted record.
$ws.Rows.Item(876).Insert()

# Populate the freshly inserted row with the new record. The date column
# is stored as literal text ("2026/02/26"), not a real date serial, so
# force the cell to Text format first to stop Excel from auto-converting
# the string into a date value, then reset the style back to Normal so the
# row carries no explicit formatting (matching every other data row).
$ws.Range("A876").NumberFormat = "@"
$ws.Range("A876").Value = "2026/02/26"
$ws.Range("A876").Style = "Normal"
$ws.Range("B876").Value = "木"
$ws.Range("C876").Value = 6
$ws.Range("D876").Value = 201
